$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.33
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 1.75
$ws.Range("J2").Value = 4.75
$ws.Range("L2").Value = 2.4
$ws.Range("X2").Value = 21
$ws.Range("Z2").Value = 41
$ws.Range("AA2").Value = 34
$ws.Range("AC2").Value = 11
$ws.Range("AI2").Value = 8.5
$ws.Range("AK2").Value = 15
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 23
$ws.Range("AS2").Value = 201
$ws.Range("AU2").Value = 8
$ws.Range("AX2").Value = 9.5
$ws.Range("H3").Value = 6.5
$ws.Range("I3").Value = 11
$ws.Range("L3").Value = 8.5
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 21
$ws.Range("O3").Value = 1.11
$ws.Range("P3").Value = 6.5
$ws.Range("Q3").Value = 1.4
$ws.Range("R3").Value = 2.88
$ws.Range("S3").Value = 1.22
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 8
$ws.Range("AA3").Value = 10
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 21
$ws.Range("AD3").Value = 13
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 251
$ws.Range("AH3").Value = 34
$ws.Range("AJ3").Value = 29
$ws.Range("AK3").Value = 126
$ws.Range("AL3").Value = 67
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 3.4
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 12
$ws.Range("AR3").Value = 29
$ws.Range("AT3").Value = 4
$ws.Range("AU3").Value = 9.5
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 251
$ws.Range("G4").Value = 2.05
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 2.75
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.8
$ws.Range("AB4").Value = 29
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 17
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 1.95
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 6.5
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 23
$ws.Range("AA5").Value = 23
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 7
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 501
$ws.Range("AH5").Value = 7.5
$ws.Range("AI5").Value = 13
$ws.Range("AL5").Value = 29
$ws.Range("AO5").Value = 15
$ws.Range("AP5").Value = 29
$ws.Range("AQ5").Value = 51
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 67
$ws.Range("AW5").Value = 4.75
$ws.Range("AY5").Value = 34
$ws.Range("AZ5").Value = 67
$ws.Range("BA5").Value = 101
$ws.Range("BB5").Value = 301
$ws.Range("G6").Value = 2.25
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 2.88
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.85
$ws.Range("R6").Value = 2
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 9
$ws.Range("AA6").Value = 17
$ws.Range("AI6").Value = 17
$ws.Range("AO6").Value = 12
$ws.Range("BD6").Value = 151
$ws.Range("G7").Value = 2.8
$ws.Range("I7").Value = 2.55
$ws.Range("L7").Value = 3.5
$ws.Range("X7").Value = 12
$ws.Range("AA7").Value = 26
$ws.Range("AJ7").Value = 11
$ws.Range("AO7").Value = 17
$ws.Range("AP7").Value = 29
$ws.Range("AQ7").Value = 51
$ws.Range("AW7").Value = 4.5
$ws.Range("AX7").Value = 17
$ws.Range("G8").Value = 1.95
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 3.7
$ws.Range("J8").Value = 2.63
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 4.5
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("AH8").Value = 10
$ws.Range("AI8").Value = 19
$ws.Range("AJ8").Value = 13
$ws.Range("AL8").Value = 34
$ws.Range("AR8").Value = 51
$ws.Range("AY8").Value = 34
$ws.Range("AZ8").Value = 81
$ws.Range("G9").Value = 2.45
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 2.9
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 3.6
$ws.Range("X9").Value = 11
$ws.Range("Y9").Value = 10
$ws.Range("Z9").Value = 23
$ws.Range("AA9").Value = 21
$ws.Range("AB9").Value = 34
$ws.Range("AC9").Value = 8.5
$ws.Range("AD9").Value = 6
$ws.Range("AH9").Value = 8.5
$ws.Range("AI9").Value = 13
$ws.Range("AJ9").Value = 11
$ws.Range("AK9").Value = 29
$ws.Range("AL9").Value = 26
$ws.Range("AO9").Value = 15
$ws.Range("AP9").Value = 26
$ws.Range("AQ9").Value = 51
$ws.Range("AS9").Value = 201
$ws.Range("AW9").Value = 4.75
$ws.Range("AX9").Value = 17
$ws.Range("AZ9").Value = 51
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 4.5
$ws.Range("N10").Value = 12
$ws.Range("Q10").Value = 1.73
$ws.Range("R10").Value = 2.08
$ws.Range("AK10").Value = 41
$ws.Range("AM10").Value = 34
$ws.Range("AN10").Value = 4
$ws.Range("AR10").Value = 51
$ws.Range("AX10").Value = 21
$ws.Range("G11").Value = 1.57
$ws.Range("H11").Value = 3.55
$ws.Range("I11").Value = 5.8
$ws.Range("J11").Value = 2.05
$ws.Range("K11").Value = 2.2
$ws.Range("L11").Value = 5.7
$ws.Range("M11").Value = 1.02
$ws.Range("N11").Value = 7.1
$ws.Range("O11").Value = 1.3
$ws.Range("P11").Value = 2.92
$ws.Range("Q11").Value = 1.88
$ws.Range("S11").Value = 1.42
$ws.Range("T11").Value = 2.47
$ws.Range("U11").Value = 1.9
$ws.Range("W11").Value = 6
$ws.Range("X11").Value = 6.9
$ws.Range("Z11").Value = 11.25
$ws.Range("AA11").Value = 13
$ws.Range("AB11").Value = 29
$ws.Range("AD11").Value = 7.1
$ws.Range("AE11").Value = 18
$ws.Range("AH11").Value = 14
$ws.Range("AI11").Value = 35
$ws.Range("AJ11").Value = 18.5
$ws.Range("AL11").Value = 70
$ws.Range("AN11").Value = 3.35
$ws.Range("AO11").Value = 7.2
$ws.Range("AQ11").Value = 22
$ws.Range("AT11").Value = 2.65
$ws.Range("AW11").Value = 7.3
$ws.Range("AX11").Value = 35
